# This script reproduces the commit "added more words into the dataset" for
# translate/record of results.xlsx
#
# Summary of the change:
#  1. The lone extra cell E202 (which held the stray shared string "s") is removed.
#  2. A whole new results block (rows 211-232) is appended right after the
#     existing block that ends at row 210. It mirrors the layout of the
#     block found at rows 189-210 (blank separator row, header row, then
#     20 data rows) but with a new set of words/accuracy numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the stray cell E202 (shared string "s") that is no longer used
# ---------------------------------------------------------------------
$ws.Range("E202").ClearContents()

# ---------------------------------------------------------------------
# 2) Duplicate the formatting (styles, borders, fonts) of the previous
#    block (rows 190-210) down onto the new block (rows 212-232) so the
#    alternating row shading / thick-bottom header-and-footer borders
#    match the existing blocks exactly.
# ---------------------------------------------------------------------
$ws.Range("A190:D210").Copy()
$ws.Range("A212:D232").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Fill in the new header + data rows (211 is left blank, just like the
#    separator row 189 before it).
# ---------------------------------------------------------------------

# Header row (same labels as every other block's header)
$ws.Range("A212").Value = "sign"
$ws.Range("B212").Value = "predicted"
$ws.Range("C212").Value = "accuracy of predicted"
$ws.Range("D212").Value = "accuracy of actual sign"

# Data rows
$data = @(
    @{ Row = 213; A = "chair";    B = "chair"; C = 1;      D = $null },
    @{ Row = 214; A = "black";    B = "black"; C = 0.9365;  D = $null },
    @{ Row = 215; A = "blue";     B = "blue";  C = 0.9492;  D = $null },
    @{ Row = 216; A = "can";      B = "can";   C = 0.9974;  D = $null },
    @{ Row = 217; A = "apple";    B = "apple"; C = 1;      D = $null },
    @{ Row = 218; A = "brown";    B = "brown"; C = 1;      D = $null },
    @{ Row = 219; A = "cat ";     B = "cat";   C = 1;      D = $null },
    @{ Row = 220; A = "cold";     B = "cold";  C = 1;      D = $null },
    @{ Row = 221; A = "come";     B = "come";  C = 1;      D = $null },
    @{ Row = 222; A = "child";    B = "child"; C = 1;      D = $null },
    @{ Row = 223; A = "Drink";    B = "apple"; C = 0.325;   D = 0.005 },
    @{ Row = 224; A = "Deaf";     B = "Deaf";  C = 1;      D = $null },
    @{ Row = 225; A = "Dog";      B = "Dog";   C = 0.99;    D = $null },
    @{ Row = 226; A = "Cow";      B = "cold";  C = 0.9984;  D = 0 },
    @{ Row = 227; A = "Eat";      B = "Eat";   C = 1;      D = $null },
    @{ Row = 228; A = "Cry";      B = "Cry";   C = 1;      D = $null },
    @{ Row = 229; A = "Drive";    B = "Drive"; C = 0.85;    D = $null },
    @{ Row = 230; A = "Egg";      B = "chair"; C = 0.9978;  D = 0.001 },
    @{ Row = 231; A = "Cup";      B = "Cup";   C = 1;      D = $null },
    @{ Row = 232; A = "Computer"; B = "cat";   C = 0.8334;  D = 0.0526 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    if ($item.D -ne $null) {
        $ws.Range("D$r").Value = $item.D
    } else {
        # No D value for this row in the source data - the pasted format
        # left an empty (but cell-less) placeholder; make sure no stray
        # D cell is written out.
        $ws.Range("D$r").ClearContents()
    }
}
